$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for a plain/unstyled cell (used to strip the quote-prefix
# marker Excel applies when a text value looks like a date, so the cell's
# style stays identical to the original unstyled date cells).
$cleanStyle = $ws.Range("A2").Style

# Row 3
$ws.Range("A3").Value = "'28-07-2022"
$ws.Range("A3").Style = $cleanStyle
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4
$ws.Range("A4").Value = "'01-08-2022"
$ws.Range("A4").Style = $cleanStyle
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("A5").Value = "'04-08-2022"
$ws.Range("A5").Style = $cleanStyle
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("A6").Value = "'08-08-2022"
$ws.Range("A6").Style = $cleanStyle

# Row 7
$ws.Range("A7").Value = "'11-08-2022"
$ws.Range("A7").Style = $cleanStyle

# Row 8
$ws.Range("A8").Value = "'15-08-2022"
$ws.Range("A8").Style = $cleanStyle

# Row 9
$ws.Range("A9").Value = "'18-08-2022"
$ws.Range("A9").Style = $cleanStyle

# Row 10
$ws.Range("A10").Value = "'22-08-2022"
$ws.Range("A10").Style = $cleanStyle

# Row 11
$ws.Range("A11").Value = "'25-08-2022"
$ws.Range("A11").Style = $cleanStyle

# Row 12
$ws.Range("A12").Value = "'29-08-2022"
$ws.Range("A12").Style = $cleanStyle

# Row 13
$ws.Range("A13").Value = "'01-09-2022"
$ws.Range("A13").Style = $cleanStyle

# Row 14
$ws.Range("A14").Value = "'05-09-2022"
$ws.Range("A14").Style = $cleanStyle

# Row 15
$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("A15").Style = $cleanStyle

# Row 16
$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("A16").Style = $cleanStyle

# Row 17
$ws.Range("A17").Value = "'15-09-2022"
$ws.Range("A17").Style = $cleanStyle

# Row 18
$ws.Range("A18").Value = "'19-09-2022"
$ws.Range("A18").Style = $cleanStyle

# Row 19
$ws.Range("A19").Value = "'22-09-2022"
$ws.Range("A19").Style = $cleanStyle

# Row 20
$ws.Range("A20").Value = "'26-09-2022"
$ws.Range("A20").Style = $cleanStyle

# Row 21
$ws.Range("A21").Value = "'29-09-2022"
$ws.Range("A21").Style = $cleanStyle
